# Apply the target edit:
# 1) Within each year group of 4 rows (A/B/C/D), swap the "B" and "C" rows,
#    i.e. swap rows 3<->4, 7<->8, 11<->12, 15<->16 (columns A:E).
# 2) Remove the extra columns F (产销率) and G (销售量) entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$swapPairs = @(
    @(3, 4),
    @(7, 8),
    @(11, 12),
    @(15, 16)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = $ws.Range("A$r1`:E$r1")
    $range2 = $ws.Range("A$r2`:E$r2")

    $tmp = $range1.Value()
    $range1.Value = $range2.Value()
    $range2.Value = $tmp
}

# Delete columns F:G completely (data + header), shrinking the used range.
$ws.Range("F1:G17").EntireColumn.Delete()
